# Update countries & provincias Spain
# - Reorders several country-name labels in the shared-string table (which,
#   since the row->index mapping is unchanged, shows up as certain rows now
#   displaying a different country name than before).
# - Refreshes the COVID-19 statistics (Casos totales, Nuevos casos, Casos
#   activos, Recuperados, Casos criticos, Muertes hoy, Muertes) for the rows
#   whose data was updated in the newer snapshot.
# - Updates the "last updated" timestamp banner in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 24 de Octubre de 2020 a las 09:54"

# Row-level updates: column A is the country label, B..H are the stats
# (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos,
# Muertes hoy, Muertes).
$updates = @(
    @{ Row = 4;   Country = "Estados Unidos";          B = 8747864; C = 911;   D = 5698965; E = 2819607; F = 0; G = 8;   H = 229292 }
    @{ Row = 7;   Country = "Rusia";                    B = 1497167; C = 16521; D = 1130818; E = 340528;  F = 0; G = 296; H = 25821 }
    @{ Row = 59;  Country = "Armenia";                  B = 75523;   C = 2213;  D = 50701;   E = 23665;   F = 0; G = 12;  H = 1157 }
    @{ Row = 60;  Country = "Austria";                  B = 74415;   C = 0;     D = 55195;   E = 18266;   F = 0; G = 0;   H = 954 }
    @{ Row = 67;  Country = "Hungria";                  B = 56098;   C = 1820;  D = 16007;   E = 38701;   F = 0; G = 38;  H = 1390 }
    @{ Row = 68;  Country = "Argelia";                  B = 55630;   C = 0;     D = 38788;   E = 14945;   F = 0; G = 0;   H = 1897 }
    @{ Row = 69;  Country = "Irlanda";                  B = 55261;   C = 0;     D = 23364;   E = 30019;   F = 0; G = 0;   H = 1878 }
    @{ Row = 70;  Country = "Kirguistan";                B = 54588;  C = 0;     D = 47050;   E = 6412;    F = 0; G = 0;   H = 1126 }
    @{ Row = 79;  Country = "Afganistan";                B = 40768;  C = 81;    D = 34023;   E = 5234;    F = 0; G = 4;   H = 1511 }
    @{ Row = 85;  Country = "El Salvador";               B = 32585;  C = 164;   D = 28258;   E = 3383;    F = 0; G = 4;   H = 944 }
    @{ Row = 88;  Country = "Australia";                 B = 27499;  C = 15;    D = 25181;   E = 1413;    F = 0; G = 0;   H = 905 }
    @{ Row = 89;  Country = "Georgia";                   B = 26503;  C = 1941;  D = 10163;   E = 16147;   F = 0; G = 10;  H = 193 }
    @{ Row = 90;  Country = "Republica de Macedonia";    B = 25991;  C = 0;     D = 18247;   E = 6861;    F = 0; G = 0;   H = 883 }
    @{ Row = 91;  Country = "Corea del Sur";             B = 25775;  C = 77;    D = 23834;   E = 1484;    F = 0; G = 2;   H = 457 }
    @{ Row = 186; Country = "Mongolia";                  B = 337;    C = 9;     D = 312;     E = 25;      F = 0; G = 0;   H = 0 }
    @{ Row = 187; Country = "Butan";                     B = 336;    C = 0;     D = 306;     E = 30;      F = 0; G = 0;   H = 0 }
    @{ Row = 216; Country = "Islas Malvinas";            B = 13;     C = 0;     D = 13;      E = 0;       F = 0; G = 0;   H = 0 }
    @{ Row = 217; Country = "Montserrat";                B = 13;     C = 0;     D = 12;      E = 0;       F = 0; G = 0;   H = 1 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).Value = $u.Country
    $ws.Cells.Item($r, 2).Value = $u.B
    $ws.Cells.Item($r, 3).Value = $u.C
    $ws.Cells.Item($r, 4).Value = $u.D
    $ws.Cells.Item($r, 5).Value = $u.E
    $ws.Cells.Item($r, 6).Value = $u.F
    $ws.Cells.Item($r, 7).Value = $u.G
    $ws.Cells.Item($r, 8).Value = $u.H
}
